# Insert a new weekly price-report row at row 454 (pushing existing rows
# 454:550 down to 455:551) for Terminal Hortofrutícola Agro Chillán - Repollo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 454 downward by inserting a fresh blank row at 454.
$ws.Rows.Item(454).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A454").Value = 7
$ws.Range("B454").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C454").Value = "Ñuble"
$ws.Range("D454").Value = 45244
$ws.Range("E454").Value = 16
$ws.Range("F454").Value = 100112006
$ws.Range("G454").Value = "Repollo"
$ws.Range("H454").Value = "Crespo record"
$ws.Range("I454").Value = "Primera"
$ws.Range("J454").Value = 400
$ws.Range("K454").Value = 1200
$ws.Range("L454").Value = 1200
$ws.Range("M454").Value = 1200
$ws.Range("N454").Value = "`$/unidad"
$ws.Range("O454").Value = "Región del Maule"
$ws.Range("P454").Value = 1200
$ws.Range("Q454").Value = 1
$ws.Range("R454").Value = "Hortaliza"
